$wb = $excel.ActiveWorkbook

# --- Rename sheets (new timestamped task-order identifiers) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16512555409842997"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-1651255543420257"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-165125554342726"

$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-1651255543485259"

$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16512555435632565"

# --- Sheet1 (GNG) stimulus file names ---
$ws1.Range("B2").Value = "go_stims-16512555409483.csv"
$ws1.Range("B3").Value = "GNG_stims-16512555409673.csv"
$ws1.Range("B4").Value = "go_stims-16512555409693.csv"
$ws1.Range("B5").Value = "GNG_stims-16512555409832993.csv"

# --- Sheet2 (NB) stimulus file names ---
$ws2.Range("B2").Value = "TB-16512555433572595.csv"
$ws2.Range("B3").Value = "ZB-match_0-16512555410492992.csv"
$ws2.Range("B4").Value = "OB-16512555419051542.csv"
$ws2.Range("B5").Value = "ZB-match_4-16512555411383002.csv"
$ws2.Range("B6").Value = "TB-16512555434012582.csv"
$ws2.Range("B7").Value = "OB-16512555418034728.csv"
$ws2.Range("B8").Value = "ZB-match_1-16512555410292997.csv"
$ws2.Range("B9").Value = "TB-16512555429062564.csv"
$ws2.Range("B10").Value = "OB-16512555416654696.csv"

# --- Sheet4 (TOL) stimulus file names ---
$ws4.Range("B2").Value = "MM_stims-16512555434522576.csv"
$ws4.Range("B3").Value = "ZM_stims-16512555434302685.csv"
$ws4.Range("B4").Value = "MM_stims-16512555434682567.csv"
$ws4.Range("B5").Value = "ZM_stims-16512555434532619.csv"
$ws4.Range("B6").Value = "MM_stims-16512555434842584.csv"
$ws4.Range("B7").Value = "ZM_stims-16512555434692585.csv"

# --- Sheet5 (vSAT) stimulus file names ---
$ws5.Range("B2").Value = "vSAT_stims-16512555435482645.csv"
$ws5.Range("B3").Value = "SAT_stims-1651255543515258.csv"
$ws5.Range("B4").Value = "vSAT_stims-16512555435312572.csv"
$ws5.Range("B5").Value = "SAT_stims-1651255543491261.csv"
